# Update resume to include Vue
# Insert "Vue.js, " into the Skills list, right after "Node.js, " and
# before "Swift,".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Node.js, Swift",        # FindText
    $true,                   # MatchCase
    $true,                   # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    "Node.js, Vue.js, Swift",# ReplaceWith
    2                        # Replace (wdReplaceAll)
)
